$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# values like "245.91" or "1.000" are not coerced into numbers.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = '30.962.97'
$ws.Cells.Item(2, 5).Value = '  +2.99%  '
$ws.Cells.Item(3, 4).Value = '1.910.62'
$ws.Cells.Item(3, 5).Value = '  +1.45%  '
$ws.Cells.Item(4, 5).Value = '  +0.39%  '
$ws.Cells.Item(5, 4).Value = '245.91'
$ws.Cells.Item(5, 5).Value = '  +0.84%  '
$ws.Cells.Item(6, 5).Value = '  +0.30%  '
$ws.Cells.Item(7, 4).Value = '0.4996'
$ws.Cells.Item(7, 5).Value = '  +0.65%  '
$ws.Cells.Item(8, 4).Value = '0.2995'
$ws.Cells.Item(8, 5).Value = '  +2.75%  '
$ws.Cells.Item(9, 4).Value = '0.06877'
$ws.Cells.Item(9, 5).Value = '  +3.91%  '
$ws.Cells.Item(10, 4).Value = '1.909.35'
$ws.Cells.Item(10, 5).Value = '  +1.43%  '
$ws.Cells.Item(11, 4).Value = '17.03'
$ws.Cells.Item(11, 5).Value = '  +1.01%  '
$ws.Cells.Item(12, 4).Value = '0.07306'
$ws.Cells.Item(12, 5).Value = '  +1.48%  '
$ws.Cells.Item(13, 4).Value = '91.47'
$ws.Cells.Item(13, 5).Value = '  +6.57%  '
$ws.Cells.Item(14, 4).Value = '5.098'
$ws.Cells.Item(14, 5).Value = '  +5.18%  '
$ws.Cells.Item(15, 4).Value = '0.6819'
$ws.Cells.Item(15, 5).Value = '  +2.36%  '
$ws.Cells.Item(16, 4).Value = '30.937.26'
$ws.Cells.Item(16, 5).Value = '  +3.01%  '
$ws.Cells.Item(17, 4).Value = '0.000008047'
$ws.Cells.Item(17, 5).Value = '  +3.37%  '
$ws.Cells.Item(18, 5).Value = '  +0.32%  '
$ws.Cells.Item(19, 5).Value = '  +3.68%  '
$ws.Cells.Item(20, 4).Value = '2.157.93'
$ws.Cells.Item(20, 5).Value = '  +1.75%  '
$ws.Cells.Item(21, 4).Value = '1.000'
$ws.Cells.Item(21, 5).Value = '  +0.31%  '
$ws.Cells.Item(22, 4).Value = '4.885'
$ws.Cells.Item(22, 5).Value = '  +2.54%  '
$ws.Cells.Item(23, 4).Value = '182.09'
$ws.Cells.Item(23, 5).Value = '  +33.83%  '
$ws.Cells.Item(24, 4).Value = '6.131'
$ws.Cells.Item(24, 5).Value = '  +9.52%  '
$ws.Cells.Item(25, 4).Value = '9.375'
$ws.Cells.Item(25, 5).Value = '  +2.48%  '
$ws.Cells.Item(26, 4).Value = '154.85'
$ws.Cells.Item(26, 5).Value = '  +3.15%  '
$ws.Cells.Item(27, 4).Value = '19.06'
$ws.Cells.Item(27, 5).Value = '  +13.60%  '
$ws.Cells.Item(28, 4).Value = '1.945'
$ws.Cells.Item(28, 5).Value = '  +2.05%  '
$ws.Cells.Item(29, 4).Value = '1.394'
$ws.Cells.Item(29, 5).Value = '  +0.94%  '
$ws.Cells.Item(30, 4).Value = '4.350'
$ws.Cells.Item(30, 5).Value = '  +4.44%  '
$ws.Cells.Item(31, 4).Value = '0.08976'
$ws.Cells.Item(31, 5).Value = '  +3.44%  '
$ws.Cells.Item(32, 4).Value = '4.056'
$ws.Cells.Item(32, 5).Value = '  +2.54%  '
$ws.Cells.Item(33, 4).Value = '0.05260'
$ws.Cells.Item(33, 5).Value = '  +5.35%  '
$ws.Cells.Item(34, 4).Value = '0.7507'
$ws.Cells.Item(34, 5).Value = '  +6.93%  '
$ws.Cells.Item(35, 5).Value = '  +3.16%  '
$ws.Cells.Item(36, 4).Value = '2.668'
$ws.Cells.Item(36, 5).Value = '  +0.54%  '
$ws.Cells.Item(37, 4).Value = '0.01976'
$ws.Cells.Item(37, 5).Value = '  +20.28%  '
$ws.Cells.Item(38, 4).Value = '2.737'
$ws.Cells.Item(38, 5).Value = '  +1.66%  '
$ws.Cells.Item(39, 4).Value = '2.180'
$ws.Cells.Item(39, 5).Value = '  -0.68%  '
$ws.Cells.Item(40, 4).Value = '0.9385'
$ws.Cells.Item(40, 5).Value = '  +0.41%  '
$ws.Cells.Item(41, 4).Value = '0.4387'
$ws.Cells.Item(41, 5).Value = '  +4.85%  '
$ws.Cells.Item(42, 4).Value = '106.57'
$ws.Cells.Item(42, 5).Value = '  +5.31%  '
$ws.Cells.Item(43, 4).Value = '5.881'
$ws.Cells.Item(43, 5).Value = '  -1.51%  '
$ws.Cells.Item(44, 5).Value = '  +0.17%  '
$ws.Cells.Item(45, 4).Value = '7.813'
$ws.Cells.Item(45, 5).Value = '  +3.64%  '
$ws.Cells.Item(46, 5).Value = '  +6.64%  '
$ws.Cells.Item(47, 4).Value = '0.05850'
$ws.Cells.Item(47, 5).Value = '  +2.33%  '
$ws.Cells.Item(48, 4).Value = '8.616'
$ws.Cells.Item(48, 5).Value = '  +5.10%  '
$ws.Cells.Item(49, 4).Value = '0.3917'
$ws.Cells.Item(49, 5).Value = '  +5.74%  '
$ws.Cells.Item(50, 4).Value = '33.29'
$ws.Cells.Item(50, 5).Value = '  +3.02%  '
$ws.Cells.Item(51, 4).Value = '1.396'
$ws.Cells.Item(51, 5).Value = '  +4.43%  '
